$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,@('Website Name', 'Profile URL')
  ,@('35photo', 'https://35photo.pro/@gorbash/')
  ,@('Audiojungle', 'https://audiojungle.net/user/gorbash')
  ,@('ask.fm', 'https://ask.fm/gorbash')
  ,@('BodyBuilding.com', 'http://api.bodybuilding.com/api-proxy/bbc/get?slug=gorbash')
  ,@('Chess.com', 'https://api.chess.com/pub/player/gorbash')
  ,@('cfx.re', 'https://forum.cfx.re/u/gorbash.json')
  ,@('chatango.com', 'https://gorbash.chatango.com')
  ,@('chaturbate', 'https://chaturbate.com/gorbash/')
  ,@('cda.pl', 'https://www.cda.pl/gorbash')
  ,@('Code Project', 'https://www.codeproject.com/Members/gorbash')
  ,@('Codewars', 'https://www.codewars.com/users/gorbash')
  ,@('cHEEZburger', 'https://profile.cheezburger.com/gorbash')
  ,@('Chyoa', 'https://chyoa.com/user/gorbash')
  ,@('Chomikuj.pl', 'https://chomikuj.pl/gorbash/')
  ,@('crowdin', 'https://crowdin.com/profile/gorbash')
  ,@('Mastodon-Defcon', 'https://defcon.social/@gorbash')
  ,@('diigo', 'https://www.diigo.com/interact_api/load_profile_info?name=gorbash')
  ,@('Discogs', 'https://api.discogs.com/users/gorbash')
  ,@('Dating.ru', 'https://dating.ru/gorbash/')
  ,@('Disqus', 'https://disqus.com/by/gorbash/')
  ,@('DockerHub', 'https://hub.docker.com/v2/users/gorbash/')
  ,@('DeviantArt', 'https://www.deviantart.com/gorbash')
  ,@('Cults3D', 'https://cults3d.com/en/users/gorbash/creations')
  ,@('datezone', 'https://www.datezone.com/users/gorbash/')
  ,@('Duolingo', 'https://www.duolingo.com/2017-06-30/users?username=gorbash&_=1628308619574')
  ,@('Demotywatory', 'https://demotywatory.pl/user/gorbash')
  ,@('eBay', 'https://www.ebay.com/usr/gorbash')
  ,@('Etsy', 'https://www.etsy.com/people/gorbash')
  ,@('Fark', 'https://www.fark.com/users/gorbash')
  ,@('FatSecret', 'https://www.fatsecret.com/member/gorbash')
  ,@('fansly', 'https://apiv2.fansly.com/api/v1/account?usernames=gorbash')
  ,@('Fabswingers', 'https://www.fabswingers.com/profile/gorbash')
  ,@('Flipboard', 'https://flipboard.com/@gorbash')
  ,@('FurAffinity', 'https://www.furaffinity.net/user/gorbash')
  ,@('Foursquare', 'https://foursquare.com/gorbash')
  ,@('Filmweb', 'https://www.filmweb.pl/user/gorbash')
  ,@('Geocaching', 'https://www.geocaching.com/p/?u=gorbash')
  ,@('Gravatar', 'https://en.gravatar.com/gorbash.json')
  ,@('GitHub', 'https://github.com/gorbash')
  ,@('GitLab', 'https://gitlab.com/gorbash')
  ,@('HackerOne', 'https://hackerone.com/gorbash')
  ,@('Hacker News', 'https://news.ycombinator.com/user?id=gorbash')
  ,@('giters', 'https://giters.com/gorbash')
  ,@('FriendFinder-X', 'https://www.friendfinder-x.com/profile/gorbash')
  ,@('IFTTT', 'https://ifttt.com/p/gorbash')
  ,@('hugging_face', 'https://huggingface.co/gorbash')
  ,@('imgur', 'https://api.imgur.com/account/v1/accounts/gorbash?client_id=546c25a59c58ad7&include=trophies%2Cmedallions')
  ,@('inaturalist', 'https://inaturalist.nz/people/gorbash')
  ,@('InkBunny', 'https://inkbunny.net/gorbash')
  ,@('issuu', 'https://issuu.com/gorbash')
  ,@('HudsonRock', 'https://cavalier.hudsonrock.com/api/json/v2/osint-tools/search-by-username?username=gorbash')
  ,@('Instagram2', 'https://dumpoir.com/v/gorbash')
  ,@('kaggle', 'https://www.kaggle.com/gorbash')
  ,@('Keybase', 'https://keybase.io/gorbash')
  ,@('Internet Archive User Search', 'https://archive.org/advancedsearch.php?q=gorbash&output=json')
  ,@('Jeuxvideo', 'https://www.jeuxvideo.com/profil/gorbash?mode=infos')
  ,@('Kongregate', 'https://www.kongregate.com/accounts/gorbash')
  ,@('Linktree', 'https://linktr.ee/gorbash')
  ,@('LibraryThing', 'https://www.librarything.com/profile/gorbash')
  ,@('Mastodon API', 'https://mastodon.social/api/v2/search?q=gorbash&limit=1&type=accounts')
  ,@('waytohey', 'https://waytohey.com/gorbash')
  ,@('MCUUID (Minecraft)', 'https://playerdb.co/api/player/minecraft/gorbash')
  ,@('Livejournal', 'https://gorbash.livejournal.com')
  ,@('Mod DB', 'https://www.moddb.com/members/gorbash')
  ,@('MyAnimeList', 'https://myanimelist.net/profile/gorbash')
  ,@('MySpace', 'https://myspace.com/gorbash')
  ,@('OpenStreetMap', 'https://www.openstreetmap.org/user/gorbash')
  ,@('Playstation Network', 'https://psnprofiles.com/xhr/search/users?q=gorbash')
  ,@('pikabu', 'https://pikabu.ru/@gorbash')
  ,@('Pinterest', 'https://www.pinterest.com/gorbash/')
  ,@('MCName (Minecraft)', 'https://mcname.info/en/search?q=gorbash')
  ,@('Periscope', 'https://www.periscope.tv/gorbash')
  ,@('Pokerstrategy', 'http://www.pokerstrategy.net/user/gorbash/profile/')
  ,@('palnet', 'https://www.palnet.io/@gorbash/')
  ,@('prv.pl', 'https://www.prv.pl/osoba/gorbash')
  ,@('public', 'https://public.com/@gorbash')
  ,@('Pornhub Users', 'https://www.pornhub.com/users/gorbash')
  ,@('Roblox', 'https://auth.roblox.com/v1/usernames/validate?username=gorbash&birthday=2019-12-31T23:00:00.000Z')
  ,@('rsi', 'https://robertsspaceindustries.com/citizens/gorbash')
  ,@('sentimente', 'https://www.sentimente.com/amp/gorbash.html')
  ,@('SFD', 'https://www.sfd.pl/profile/gorbash')
  ,@('RumbleUser', 'https://rumble.com/user/gorbash')
  ,@('Snapchat Stories', 'https://story.snapchat.com/s/gorbash')
  ,@('sofurry', 'https://gorbash.sofurry.com')
  ,@('slideshare', 'https://www.slideshare.net/gorbash')
  ,@('Sourceforge', 'https://sourceforge.net/u/gorbash/profile')
  ,@('Steam', 'https://steamcommunity.com/id/gorbash')
  ,@('SoundCloud', 'https://soundcloud.com/gorbash')
  ,@('Spotify', 'https://open.spotify.com/user/gorbash')
  ,@('tabletoptournament', 'https://www.tabletoptournaments.net/eu/player/gorbash')
  ,@('Telegram', 'https://t.me/gorbash')
  ,@('themeforest', 'https://themeforest.net/user/gorbash')
  ,@('TikTok', 'https://www.tiktok.com/oembed?url=https://www.tiktok.com/@gorbash')
  ,@('Trello', 'https://trello.com/1/Members/gorbash?fields=activityBlocked%2CavatarUrl%2Cbio%2CbioData%2Cconfirmed%2CfullName%2CidEnterprise%2CidMemberReferrer%2Cinitials%2CmemberType%2CnonPublic%2Cproducts%2Curl%2Cusername')
  ,@('trakt', 'https://trakt.tv/users/gorbash')
  ,@('tumblr', 'https://gorbash.tumblr.com')
  ,@('Twitter archived profile', 'http://archive.org/wayback/available?url=https://twitter.com/gorbash')
  ,@('Twitter archived tweets', 'http://archive.org/wayback/available?url=https://twitter.com/gorbash/status/*')
  ,@('Venmo', 'https://account.venmo.com/u/gorbash')
  ,@('untappd', 'https://untappd.com/user/gorbash/')
  ,@('Vimeo', 'https://vimeo.com/gorbash')
  ,@('Twitter', 'https://nitter.privacydev.net/gorbash')
  ,@('Twitch', 'https://twitchtracker.com/gorbash')
  ,@('wattpad', 'https://www.wattpad.com/api/v3/users/gorbash?fields=username%2Cname%2Cdescription%2Cavatar%2CbackgroundUrl%2CcreateDate%2Clocation%2Cfollowing%2CfollowingRequest%2CnumFollowing%2Cfollower%2CfollowerRequest%2CnumFollowers%2CnumLists%2CnumStoriesPublished%2CvotesReceived%2Cfacebook%2Ctwitter%2Cwebsite%2Csmashwords%2Chighlight_colour%2Chtml_enabled%2Cverified%2Cambassador%2Cwattpad_squad%2Cis_staff%2Cprograms(wattpad_stars)%2CisPrivate%2CisMuted%2CexternalId%2Cnotes')
  ,@('Wikidot', 'http://www.wikidot.com/user:info/gorbash')
  ,@('WordPress Support', 'https://wordpress.org/support/users/gorbash/')
  ,@('VK', 'https://vk.com/gorbash')
  ,@('xHamster', 'https://xhamster.com/users/gorbash')
  ,@('YouTube User', 'https://www.youtube.com/user/gorbash/about')
  ,@('XVIDEOS-profiles', 'https://www.xvideos.com/profiles/gorbash')
  ,@('YouTube User2', 'https://www.youtube.com/@gorbash')
  ,@('XVIDEOS-models', 'https://www.xvideos.com/models/gorbash')
  ,@('Zbiornik', 'https://mini.zbiornik.com/gorbash')
  ,@('XNXX', 'https://www.xnxx.com/mobile/profile/gorbash')
  ,@('Xbox Gamertag', 'https://www.xboxgamertag.com/search/gorbash')
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 1).Value = $data[$i][0]
  $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Range("A117:B200").ClearContents()
